# "minor fix in slides"
#
# 1) Slides 4 and 5 swap places: the short "Consider a linear combination..."
#    slide moves up to position 4, and the longer "A linear combination..."
#    derivation slide (previously hidden) drops to position 5 and is
#    unhidden so it now shows in the deck.
# 2) The speaker notes on the (now unhidden, now-position-5) derivation
#    slide no longer carry the old "Interpretation what we are doing,
#    loadings, scores" placeholder reminder text.
# 3) The cached date field on the notes master is refreshed.

$p = $ppt.ActivePresentation

# --- 1. Reorder: move slide 5 ("Consider a linear combination...") so it
#        lands right before slide 4 ("A linear combination..." derivation).
$shortSlide = $p.Slides.Item(5)
$shortSlide.MoveTo(4)

# The derivation slide (now at position 5) was previously marked hidden;
# it should be a normal, visible slide now.
$longSlide = $p.Slides.Item(5)
$longSlide.SlideShowTransition.Hidden = 0

# --- 2. Clear the leftover "Interpretation..." reminder in that slide's
#        speaker notes now that the real content is in place.
$notes = $longSlide.NotesPage
$notesBody = $notes.Shapes.Item(2)
$notesBody.TextFrame.TextRange.Text = ""

# --- 3. Refresh the cached date shown on the notes master.
$notesMaster = $p.NotesMaster
$dateShape = $notesMaster.Shapes.Item(2)
$dateShape.TextFrame.TextRange.Text = "20-08-2022"
